$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Update existing Winner / 2nd / 3rd place results with times appended ---
$ws.Range("B2").Value = "Noah Lyles - 9.79"
$ws.Range("C2").Value = "Kishane Thompson - 9.79"
$ws.Range("D2").Value = "Fred Kerley - 9.81"

$ws.Range("B3").Value = "Julien Alfred - 10.72"
$ws.Range("C3").Value = "Sha'Carri Richardson - 10.87"
$ws.Range("D3").Value = "Melissa Jefferson - 10.92"

$ws.Range("B5").Value = "Gabrielle Thomas - 21.83"
$ws.Range("C5").Value = "Julien Alfred - 22.08"
$ws.Range("D5").Value = "Brittany Brown - 22.20"

# --- New 4th-8th place headers ---
$ws.Range("E1").Value = "4th Place"
$ws.Range("F1").Value = "5th Place"
$ws.Range("G1").Value = "6th Place"
$ws.Range("H1").Value = "7th Place"
$ws.Range("I1").Value = "8th Place"

# --- 100m - M row (row 2) ---
$ws.Range("E2").Value = "Akani Simbine - 9.82"
$ws.Range("F2").Value = "Lamont Jacobs - 9.85"
$ws.Range("G2").Value = "Letsile Tebogo - 9.86"
$ws.Range("H2").Value = "Kenneth Bednarek - 9.88"
$ws.Range("I2").Value = "Oblique Seville - 9.91"

# --- 100m - F row (row 3) ---
$ws.Range("E3").Value = "Daryll Neita - 10.96"
$ws.Range("F3").Value = "Twanisga Terry - 10.97"
$ws.Range("G3").Value = "Mujinga Kambundji - 10.99"
$ws.Range("H3").Value = "Tia Clayton - 11.04"
$ws.Range("I3").Value = "Marie-Josee Lou-Smith - 13.84"

# --- 200m - F row (row 5) ---
$ws.Range("E5").Value = "Dina Asher-Smith - 22.22"
$ws.Range("F5").Value = "Daryll Neita - 22.23"
$ws.Range("G5").Value = "Favour Ofili - 22.24"
$ws.Range("H5").Value = "Mckenzie Long - 22.42"
$ws.Range("I5").Value = "Jessika Gbai - 22.70"

# --- Remove the bold/sz10 header & row-label fonts, reverting to default font ---
$ws.Range("A1:D1").Font.Size = 11
$ws.Range("A1:D1").Font.Bold = $false

$ws.Range("A2:A9").Font.Size = 11
$ws.Range("A2:A9").Font.Bold = $false

# --- B5 picks up an explicit black font color (theme -> explicit RGB) with no border (artifact of direct edit) ---
$ws.Range("B5").Font.Color = 0
$ws.Range("B5").Borders.LineStyle = -4142

# --- Column widths for the new columns ---
$ws.Columns.Item(5).ColumnWidth = 17.33
$ws.Columns.Item(6).ColumnWidth = 18.89
$ws.Columns.Item(7).ColumnWidth = 22.22
$ws.Columns.Item(8).ColumnWidth = 20.66
$ws.Columns.Item(9).ColumnWidth = 25.55
